# Apply weekly update to the Espinaca price sheet:
# - Insert a brand new observation at row 68 (most recent week), shifting the
#   existing time-series (columns D, I:P) for rows 68-192 down by one row.
# - The oldest observation (previously in row 192) is appended as a brand new
#   row 193, keeping all its original column values.
# - Columns A, B, C, E, F, G, H, Q, R are identical on every data row, so the
#   new row 193 simply reuses those constants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 68
$lastRow = 192
$newLastRow = 193

# --- 1. Capture the current ("old") values before we overwrite anything ---

# Full old last row (192), used to build the brand new row 193.
$oldLastRowVals = $ws.Range("A$lastRow`:R$lastRow").Value2

# Old D68:D191 and I68:P191 -> these become the new D69:D192 / I69:P192.
$oldD = $ws.Range("D$firstRow`:D$($lastRow-1)").Value2
$oldIP = $ws.Range("I$firstRow`:P$($lastRow-1)").Value2

# --- 2. Create the brand new row 193 from the old row 192 values ---

$ws.Range("A$newLastRow").Value = $oldLastRowVals[1,1]
$ws.Range("B$newLastRow").Value = $oldLastRowVals[1,2]
$ws.Range("C$newLastRow").Value = $oldLastRowVals[1,3]
$ws.Range("D$newLastRow").Value = $oldLastRowVals[1,4]
$ws.Range("E$newLastRow").Value = $oldLastRowVals[1,5]
$ws.Range("F$newLastRow").Value = $oldLastRowVals[1,6]
$ws.Range("G$newLastRow").Value = $oldLastRowVals[1,7]
$ws.Range("H$newLastRow").Value = $oldLastRowVals[1,8]
$ws.Range("I$newLastRow").Value = $oldLastRowVals[1,9]
$ws.Range("J$newLastRow").Value = $oldLastRowVals[1,10]
$ws.Range("K$newLastRow").Value = $oldLastRowVals[1,11]
$ws.Range("L$newLastRow").Value = $oldLastRowVals[1,12]
$ws.Range("M$newLastRow").Value = $oldLastRowVals[1,13]
$ws.Range("N$newLastRow").Value = $oldLastRowVals[1,14]
$ws.Range("O$newLastRow").Value = $oldLastRowVals[1,15]
$ws.Range("P$newLastRow").Value = $oldLastRowVals[1,16]
$ws.Range("Q$newLastRow").Value = $oldLastRowVals[1,17]
$ws.Range("R$newLastRow").Value = $oldLastRowVals[1,18]

# Match the date number format used by the rest of column D.
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat

# --- 3. Shift the time-series down by one row: new[r] = old[r-1] ---

$ws.Range("D$($firstRow+1):D$lastRow").Value = $oldD
$ws.Range("I$($firstRow+1):P$lastRow").Value = $oldIP

# --- 4. Write the brand new observation into row 68 ---

$ws.Range("D$firstRow").Value = 44469
$ws.Range("J$firstRow").Value = 120
$ws.Range("K$firstRow").Value = 3000
$ws.Range("L$firstRow").Value = 3000
$ws.Range("M$firstRow").Value = 3000
$ws.Range("P$firstRow").Value = 1000

Write-Host "Shift complete. New used range rows:" $ws.UsedRange.Rows.Count
